$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.637.95'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '3.008.25'
$ws.Range('E3').Value = '  -3.70%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.06'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.51'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.12%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.563'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.10%  '
$ws.Range('D9').Value = '3.009.42'
$ws.Range('E9').Value = '  -3.47%  '
$ws.Range('E10').Value = '  -2.17%  '
$ws.Range('E11').Value = '  -5.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.36%  '
$ws.Range('D13').Value = '3.521.59'
$ws.Range('E13').Value = '  -3.92%  '
$ws.Range('E14').Value = '  -3.36%  '
$ws.Range('D15').Value = '62.736.63'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.86'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.51%  '
$ws.Range('D17').Value = '3.002.96'
$ws.Range('E17').Value = '  -3.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000149'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '394.86'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.09'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.84'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.61'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.68%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.98'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.465'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.48%  '
$ws.Range('E26').Value = '  -6.53%  '
$ws.Range('D27').Value = '0.0₃0961'
$ws.Range('E27').Value = '  -3.81%  '
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E31').Value = '  -2.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.41'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '161.47'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.66'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.00'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('E36').Value = '  -1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.29'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.12%  '
$ws.Range('E38').Value = '  -3.85%  '
$ws.Range('D39').Value = '2.472.47'
$ws.Range('E39').Value = '  -9.67%  '
$ws.Range('E40').Value = '  -2.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.46'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.89'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.63%  '
$ws.Range('E43').Value = '  -4.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0595'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0247'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.99'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -6.77%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.68'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.82%  '
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.51'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0944'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '262.59'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.12%  '
